# Comments & MU Fields
# Added a New MU Field in Import, Export, Add, Update journeys.
# Changed Import Format.
# Ability to see comments in exported file

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Sheet1"): add two new header columns: Comments (O1) and MU (P1) ---
# Copy the formatting of the last existing header cell (N1) onto the two new
# header cells so they pick up the same style (border/fill/font/alignment).
$ws1.Range("N1").Copy() | Out-Null
$ws1.Range("O1:P1").PasteSpecial(-4122) | Out-Null
$ws1.Application.CutCopyMode = $false

$ws1.Range("O1").Value = "Comments"
$ws1.Range("P1").Value = "MU"

# --- Sheet1: view / selection ---
$ws1.Activate() | Out-Null
$ws1.Range("B1").Select() | Out-Null

# --- Sheet1: page setup ---
$ws1.PageSetup.LeftMargin = 0.1 * 72
$ws1.PageSetup.RightMargin = 0.1 * 72
$ws1.PageSetup.TopMargin = 0 * 72
$ws1.PageSetup.BottomMargin = 0.75 * 72
$ws1.PageSetup.HeaderMargin = 0.25 * 72
$ws1.PageSetup.FooterMargin = 0.3 * 72
$ws1.PageSetup.PaperSize = 9

# --- Sheet1: column widths (narrower columns in the updated import format) ---
# ColumnWidth is expressed in characters; the stored xlsx "width" adds the
# standard 5/6-character cell padding on save.
$pad = 0.8333333333333334
$ws1.Columns.Item(1).ColumnWidth = 7.46875 - $pad
$ws1.Columns.Item(2).ColumnWidth = 12.703125 - $pad
$ws1.Columns.Item(3).ColumnWidth = 10.5859375 - $pad
$ws1.Columns.Item(4).ColumnWidth = 5.8203125 - $pad
$ws1.Columns.Item(5).ColumnWidth = 32.8203125 - $pad
$ws1.Columns.Item(6).ColumnWidth = 6.9375 - $pad
$ws1.Columns.Item(7).ColumnWidth = 7.52734375 - $pad
$ws1.Columns.Item(8).ColumnWidth = 7.29296875 - $pad
$ws1.Columns.Item(9).ColumnWidth = 6.46875 - $pad
$ws1.Columns.Item(10).ColumnWidth = 7.46875 - $pad
$ws1.Columns.Item(11).ColumnWidth = 6.52734375 - $pad
$ws1.Columns.Item(12).ColumnWidth = 7.17578125 - $pad
$ws1.Columns.Item(13).ColumnWidth = 5.29296875 - $pad
$ws1.Columns.Item(14).ColumnWidth = 6.8203125 - $pad
$ws1.Columns.Item(15).ColumnWidth = 19.234375 - $pad
$ws1.Columns.Item(16).ColumnWidth = 6.234375 - $pad

# --- Sheet2 ("Style codes"): page setup gets an explicit pageSetup element ---
$ws2.PageSetup.Orientation = 1

# Minor column width tweaks on sheet2
$ws2.Columns.Item(2).ColumnWidth = 13.46875 - $pad

Write-Host "Applied Comments & MU fields update."
